$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.81824933333333
$ws.Range("H2").Value = 53.454748
$ws.Range("I2").Value = 0.05180179233147032
$ws.Range("J2").Value = 0.05180179233147034
$ws.Range("M2").Value = 53.419777
$ws.Range("N2").Value = 160.259331
$ws.Range("O2").Value = 0.9055031838237333
$ws.Range("P2").Value = 0.9055031838237334
$ws.Range("Q2").Value = 951.8469059170654
$ws.Range("R2").Value = 8566.622153253589
$ws.Range("S2").Value = 0.04690668788392223
$ws.Range("T2").Value = 0.04690668788392226
$ws.Range("G3").Value = 17.81824933333333
$ws.Range("H3").Value = 53.454748
$ws.Range("I3").Value = 0.05180179233147032
$ws.Range("J3").Value = 0.05180179233147034
$ws.Range("O3").Value = 0.02821782544644075
$ws.Range("P3").Value = 0.02821782544644075
$ws.Range("Q3").Value = 29.66201590753422
$ws.Range("R3").Value = 266.958143167808
$ws.Range("S3").Value = 0.001461733933822203
$ws.Range("T3").Value = 0.001461733933822203
$ws.Range("G4").Value = 17.81824933333333
$ws.Range("H4").Value = 53.454748
$ws.Range("I4").Value = 0.05180179233147032
$ws.Range("J4").Value = 0.05180179233147034
$ws.Range("M4").Value = 3.864071
$ws.Range("N4").Value = 11.592213
$ws.Range("O4").Value = 0.0654987495178229
$ws.Range("P4").Value = 0.06549874951782292
$ws.Range("Q4").Value = 68.85098051970267
$ws.Range("R4").Value = 619.6588246773241
$ws.Range("S4").Value = 0.003392952620493254
$ws.Range("T4").Value = 0.003392952620493256
$ws.Range("G5").Value = 17.81824933333333
$ws.Range("H5").Value = 53.454748
$ws.Range("I5").Value = 0.05180179233147032
$ws.Range("J5").Value = 0.05180179233147034
$ws.Range("M5").Value = 0.04602999999999999
$ws.Range("N5").Value = 0.13809
$ws.Range("O5").Value = 0.0007802412120029336
$ws.Range("P5").Value = 0.0007802412120029337
$ws.Range("Q5").Value = 0.8201740168133332
$ws.Range("R5").Value = 7.381566151319999
$ws.Range("S5").Value = 0.00004041789323263068
$ws.Range("T5").Value = 0.0000404178932326307
$ws.Range("I6").Value = 0.4402211587141748
$ws.Range("J6").Value = 0.4402211587141748
$ws.Range("M6").Value = 53.419777
$ws.Range("N6").Value = 160.259331
$ws.Range("O6").Value = 0.9055031838237333
$ws.Range("P6").Value = 0.9055031838237334
$ws.Range("Q6").Value = 8088.970071924521
$ws.Range("R6").Value = 72800.73064732068
$ws.Range("S6").Value = 0.3986216608022582
$ws.Range("T6").Value = 0.3986216608022584
$ws.Range("I7").Value = 0.4402211587141748
$ws.Range("J7").Value = 0.4402211587141748
$ws.Range("O7").Value = 0.02821782544644075
$ws.Range("P7").Value = 0.02821782544644075
$ws.Range("S7").Value = 0.01242208381442647
$ws.Range("T7").Value = 0.01242208381442648
$ws.Range("I8").Value = 0.4402211587141748
$ws.Range("J8").Value = 0.4402211587141748
$ws.Range("M8").Value = 3.864071
$ws.Range("N8").Value = 11.592213
$ws.Range("O8").Value = 0.0654987495178229
$ws.Range("P8").Value = 0.06549874951782292
$ws.Range("Q8").Value = 585.1082956559601
$ws.Range("R8").Value = 5265.97466090364
$ws.Range("S8").Value = 0.02883393540706549
$ws.Range("T8").Value = 0.0288339354070655
$ws.Range("I9").Value = 0.4402211587141748
$ws.Range("J9").Value = 0.4402211587141748
$ws.Range("M9").Value = 0.04602999999999999
$ws.Range("N9").Value = 0.13809
$ws.Range("O9").Value = 0.0007802412120029336
$ws.Range("P9").Value = 0.0007802412120029337
$ws.Range("Q9").Value = 6.9699896428
$ws.Range("R9").Value = 62.72990678519999
$ws.Range("S9").Value = 0.0003434786904244835
$ws.Range("T9").Value = 0.0003434786904244836
$ws.Range("G10").Value = 100.6958183333333
$ws.Range("H10").Value = 302.087455
$ws.Range("I10").Value = 0.29274614875843
$ws.Range("J10").Value = 0.2927461487584301
$ws.Range("M10").Value = 53.419777
$ws.Range("N10").Value = 160.259331
$ws.Range("O10").Value = 0.9055031838237333
$ws.Range("P10").Value = 0.9055031838237334
$ws.Range("Q10").Value = 5379.14816019918
$ws.Range("R10").Value = 48412.33344179261
$ws.Range("S10").Value = 0.2650825697528946
$ws.Range("T10").Value = 0.2650825697528947
$ws.Range("G11").Value = 100.6958183333333
$ws.Range("H11").Value = 302.087455
$ws.Range("I11").Value = 0.29274614875843
$ws.Range("J11").Value = 0.2927461487584301
$ws.Range("O11").Value = 0.02821782544644075
$ws.Range("P11").Value = 0.02821782544644075
$ws.Range("Q11").Value = 167.6281945184089
$ws.Range("R11").Value = 1508.65375066568
$ws.Range("S11").Value = 0.008260659725783156
$ws.Range("T11").Value = 0.008260659725783159
$ws.Range("G12").Value = 100.6958183333333
$ws.Range("H12").Value = 302.087455
$ws.Range("I12").Value = 0.29274614875843
$ws.Range("J12").Value = 0.2927461487584301
$ws.Range("M12").Value = 3.864071
$ws.Range("N12").Value = 11.592213
$ws.Range("O12").Value = 0.0654987495178229
$ws.Range("P12").Value = 0.06549874951782292
$ws.Range("Q12").Value = 389.0957914431018
$ws.Range("R12").Value = 3501.862122987916
$ws.Range("S12").Value = 0.01917450666983573
$ws.Range("T12").Value = 0.01917450666983574
$ws.Range("G13").Value = 100.6958183333333
$ws.Range("H13").Value = 302.087455
$ws.Range("I13").Value = 0.29274614875843
$ws.Range("J13").Value = 0.2927461487584301
$ws.Range("M13").Value = 0.04602999999999999
$ws.Range("N13").Value = 0.13809
$ws.Range("O13").Value = 0.0007802412120029336
$ws.Range("P13").Value = 0.0007802412120029337
$ws.Range("Q13").Value = 4.635028517883334
$ws.Range("R13").Value = 41.71525666095
$ws.Range("S13").Value = 0.0002284126099164686
$ws.Range("T13").Value = 0.0002284126099164686
$ws.Range("G14").Value = 74.032918
$ws.Range("H14").Value = 222.098754
$ws.Range("I14").Value = 0.2152309001959248
$ws.Range("J14").Value = 0.2152309001959249
$ws.Range("M14").Value = 53.419777
$ws.Range("N14").Value = 160.259331
$ws.Range("O14").Value = 0.9055031838237333
$ws.Range("P14").Value = 0.9055031838237334
$ws.Range("Q14").Value = 3954.821970219286
$ws.Range("R14").Value = 35593.39773197357
$ws.Range("S14").Value = 0.1948922653846581
$ws.Range("T14").Value = 0.1948922653846582
$ws.Range("G15").Value = 74.032918
$ws.Range("H15").Value = 222.098754
$ws.Range("I15").Value = 0.2152309001959248
$ws.Range("J15").Value = 0.2152309001959249
$ws.Range("O15").Value = 0.02821782544644075
$ws.Range("P15").Value = 0.02821782544644075
$ws.Range("Q15").Value = 123.2424998840426
$ws.Range("R15").Value = 1109.182498956384
$ws.Range("S15").Value = 0.006073347972408917
$ws.Range("T15").Value = 0.00607334797240892
$ws.Range("G16").Value = 74.032918
$ws.Range("H16").Value = 222.098754
$ws.Range("I16").Value = 0.2152309001959248
$ws.Range("J16").Value = 0.2152309001959249
$ws.Range("M16").Value = 3.864071
$ws.Range("N16").Value = 11.592213
$ws.Range("O16").Value = 0.0654987495178229
$ws.Range("P16").Value = 0.06549874951782292
$ws.Range("Q16").Value = 286.068451489178
$ws.Range("R16").Value = 2574.616063402602
$ws.Range("S16").Value = 0.01409735482042842
$ws.Range("T16").Value = 0.01409735482042843
$ws.Range("G17").Value = 74.032918
$ws.Range("H17").Value = 222.098754
$ws.Range("I17").Value = 0.2152309001959248
$ws.Range("J17").Value = 0.2152309001959249
$ws.Range("M17").Value = 0.04602999999999999
$ws.Range("N17").Value = 0.13809
$ws.Range("O17").Value = 0.0007802412120029336
$ws.Range("P17").Value = 0.0007802412120029337
$ws.Range("Q17").Value = 3.407735215539999
$ws.Range("R17").Value = 30.66961693986
$ws.Range("S17").Value = 0.0001679320184293509
$ws.Range("T17").Value = 0.0001679320184293509
